$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'259.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Formula = "'1.20%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Formula = "'26.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Formula = "'-0.44%"
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Formula = "'0.79%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Formula = "'0.06069"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Formula = "'3.25%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Formula = "'6.700"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Formula = "'1.19%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Formula = "'0.8600"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Formula = "'0.15%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Formula = "'0.9223"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Formula = "'-2.14%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Formula = "'-0.10%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Formula = "'0.05329"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Formula = "'23.95%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Formula = "'0.07113"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Formula = "'0.31%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Formula = "'0.03133"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Formula = "'-0.18%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Formula = "'0.09143"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Formula = "'0.001533"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Formula = "'-0.48%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Formula = "'0.0006047"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Formula = "'-0.07%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Formula = "'0.006061"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Formula = "'-2.64%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Formula = "'-1.05%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Formula = "'-1.12%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Formula = "'-1.27%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Formula = "'2.43%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Formula = "'0.1298"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Formula = "'-0.17%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Formula = "'4.093"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Formula = "'7.10%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Formula = "'0.04232"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Formula = "'-0.23%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Formula = "'0.001217"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Formula = "'-0.30%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Formula = "'0.004042"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Formula = "'-5.65%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Formula = "'0.0001200"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Formula = "'-0.09%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Formula = "'-21.37%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D40").Formula = "'0.03872"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Formula = "'1.27%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Formula = "'0.1118"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Formula = "'1.52%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Formula = "'0.004129"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Formula = "'-33.56%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Formula = "'0.01494"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Formula = "'30.90%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Formula = "'-9.57%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Formula = "'0.00005408"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Formula = "'-1.25%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Formula = "'-0.07%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Formula = "'-19.22%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Formula = "'-37.12%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Formula = "'-0.07%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Formula = "'-0.07%"
$ws.Range("E50").Style = "Normal"
